# Insert a new data row at row 462 (pushing existing rows 462..497 down to 463..498)
# and populate it with the new record's values, matching the surrounding rows'
# layout (Macroferia Regional de Talca - Repollo, Maule).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 462 downward by inserting a new blank row at position 462.
$ws.Rows("462:462").Insert()

# Fill in the new row's values.
$ws.Cells.Item(462, 1).Value = 5
$ws.Cells.Item(462, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(462, 3).Value = "Maule"
$ws.Cells.Item(462, 4).Value = 45106
$ws.Cells.Item(462, 5).Value = 7
$ws.Cells.Item(462, 6).Value = 100112006
$ws.Cells.Item(462, 7).Value = "Repollo"
$ws.Cells.Item(462, 8).Value = "Crespo record"
$ws.Cells.Item(462, 9).Value = "Primera"
$ws.Cells.Item(462, 10).Value = 5000
$ws.Cells.Item(462, 11).Value = 600
$ws.Cells.Item(462, 12).Value = 600
$ws.Cells.Item(462, 13).Value = 600
$ws.Cells.Item(462, 14).Value = "`$/unidad"
$ws.Cells.Item(462, 15).Value = "Región del Maule"
$ws.Cells.Item(462, 16).Value = 600
$ws.Cells.Item(462, 17).Value = 1
$ws.Cells.Item(462, 18).Value = "Hortaliza"
